$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("A14").Value = "2024/109609.5"
$ws.Range("B14").Value = "SECRETARIA 2-B"
$ws.Range("C14").Value = 45390
$ws.Range("C14").NumberFormat = "YYYY-MM-DD"
$ws.Range("D14").Value = "SANDERLAN"
$ws.Range("E14").Value = "Fulano1`n"

# Row 15
$ws.Range("A15").Value = "2024/029609.9"
$ws.Range("B15").Value = "SECRETARIA 2-B"
$ws.Range("C15").Value = 45390
$ws.Range("C15").NumberFormat = "YYYY-MM-DD"
$ws.Range("D15").Value = "MARCOS"
$ws.Range("E15").Value = "Fulano4`n"
